$d = $word.ActiveDocument

# Anchor on the rent amount ("1.000.") that closes out "§ 3 Miete", so the
# new "§ 4 Mietbeginn" section lands right after it no matter how the rest
# of the document is laid out.
$find = $d.Content
[void]$find.Find.Execute("1.000.", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

# The matched run lives in the document's last paragraph; that paragraph's
# Range includes its trailing paragraph mark, so its End is the true
# end-of-body insertion point (one past the last visible character).
$endOfBody = $d.Paragraphs.Last.Range.End

# Re-anchor via a brand-new, zero-length Range built from that offset.
# (Re-using/collapsing the paragraph's own Range object at end-of-body can
# make the host treat the insert as "replace this paragraph" instead of
# "append new paragraphs after it".)
$anchor = $d.Range($endOfBody, $endOfBody)

# Append: one empty paragraph, the new heading, and the new body text —
# as plain WordprocessingML so the empty paragraph stays truly empty
# (no stray run) just like the diff shows.
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newParagraphs = "<w:p $w/>" `
    + "<w:p $w><w:r><w:t>§ 4 Mietbeginn</w:t></w:r></w:p>" `
    + "<w:p $w><w:r><w:t>Das Mietverhältnis beginnt am 29. Juni 2023.</w:t></w:r></w:p>"

[void]$anchor.InsertXML($newParagraphs)
